{"js": "// Remove the \"Freixo-de-folha-estreita\", \"Amieiro\", \"Turfeira\",\n// \"Bruchia vogesiaca\" and \"Veronica micrantha\" sections (each a red,\n// centered heading paragraph followed by its body paragraph(s)) that\n// follow the \"Carvalheira\" section at the end of the document body.\n// After the edit, the document should end with the \"Carvalheira\"\n// paragraph, the trailing empty paragraph that already follows it,\n// and then the section properties.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the first paragraph to remove: the \"Freixo-de-folha-estreita\"\n// heading. Everything from this paragraph through the end of the body\n// is deleted.\nlet startIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.trim() === \"Freixo-de-folha-estreita\") {\n    startIndex = i;\n    break;\n  }\n}\n\nif (startIndex === -1) {\n  throw new Error('Could not find the \"Freixo-de-folha-estreita\" heading paragraph.');\n}\n\n// Delete paragraphs from the end of the body back up to (and\n// including) the heading paragraph, so indices remain valid while we\n// iterate.\nfor (let i = items.length - 1; i >= startIndex; i--) {\n  items[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"Freixo-de-folha-estreita\", \"Amieiro\", \"Turfeira\",\n# \"Bruchia vogesiaca\" and \"Veronica micrantha\" sections (each a red,\n# centered heading paragraph followed by its body paragraph(s)) that\n# follow the \"Carvalheira\" section at the end of the document body.\n# After the edit, the document should end with the \"Carvalheira\"\n# paragraph, the trailing empty paragraph that already follows it,\n# and then the section properties.\n\n$d = $word.ActiveDocument\n\n# Locate the \"Freixo-de-folha-estreita\" heading paragraph - the first\n# paragraph of the content block that needs to be removed.\n$searchRange = $d.Content\n$searchRange.Find.ClearFormatting()\n$searchRange.Find.Text = \"Freixo-de-folha-estreita\"\n$searchRange.Find.MatchCase = $true\n$searchRange.Find.MatchWholeWord = $false\n$searchRange.Find.Forward = $true\n$searchRange.Find.Wrap = 0\n$found = $searchRange.Find.Execute()\n\nif (-not $found) {\n    throw \"Could not find the 'Freixo-de-folha-estreita' heading paragraph.\"\n}\n\n# Build a range spanning from the start of that heading through the\n# very end of the document, and delete it in one go.\n$deleteRange = $d.Range($searchRange.Start, $d.Content.End)\n$deleteRange.Delete()\n"}
